# Auto-generated edit script: update Leve profit-calculation columns (H-N)
# across all profession sheets, per the scheduled price-refresh diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1149.6666
$ws.Range("I4").Value = 1624.5
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 1624.5
$ws.Range("L4").Value = 200
$ws.Range("M4").Value = -1510.5
$ws.Range("N4").Value = -428
$ws.Range("H6").Value = 11013.941
$ws.Range("I6").Value = 12095.733
$ws.Range("J6").Value = 2900.5
$ws.Range("K6").Value = 36287.199
$ws.Range("L6").Value = 8701.5
$ws.Range("M6").Value = -36175.199
$ws.Range("N6").Value = -8925.5
$ws.Range("H32").Value = 24570.615
$ws.Range("I32").Value = 44836.832
$ws.Range("J32").Value = 7199.5713
$ws.Range("K32").Value = 44836.832
$ws.Range("L32").Value = 7199.5713
$ws.Range("M32").Value = -44510.832
$ws.Range("N32").Value = -7851.5713
$ws.Range("H40").Value = 2032.0883
$ws.Range("I40").Value = 1990.125
$ws.Range("J40").Value = 2069.389
$ws.Range("K40").Value = 1990.125
$ws.Range("L40").Value = 2069.389
$ws.Range("M40").Value = -1815.125
$ws.Range("N40").Value = -2419.389
$ws.Range("H41").Value = 433.33334
$ws.Range("I41").Value = 522.8333
$ws.Range("K41").Value = 522.8333
$ws.Range("M41").Value = -82.83330000000001
$ws.Range("H53").Value = 463.92593
$ws.Range("I53").Value = 439.13043
$ws.Range("J53").Value = 606.5
$ws.Range("K53").Value = 439.13043
$ws.Range("L53").Value = 606.5
$ws.Range("M53").Value = 197.86957
$ws.Range("N53").Value = -1880.5
$ws.Range("H55").Value = 126.35714
$ws.Range("I55").Value = 126.916664
$ws.Range("J55").Value = 123
$ws.Range("K55").Value = 126.916664
$ws.Range("L55").Value = 123
$ws.Range("M55").Value = 87.083336
$ws.Range("N55").Value = -551
$ws.Range("H112").Value = 2336.125
$ws.Range("I112").Value = 2299.5
$ws.Range("J112").Value = 2348.3333
$ws.Range("K112").Value = 6898.5
$ws.Range("L112").Value = 7044.999899999999
$ws.Range("M112").Value = -5790.5
$ws.Range("N112").Value = -9260.999899999999
$ws.Range("H116").Value = 8101.3
$ws.Range("I116").Value = 6673.875
$ws.Range("J116").Value = 9052.916999999999
$ws.Range("K116").Value = 6673.875
$ws.Range("L116").Value = 9052.916999999999
$ws.Range("M116").Value = -3231.875
$ws.Range("N116").Value = -15936.917
$ws.Range("H132").Value = 3247
$ws.Range("I132").Value = 3412.8667
$ws.Range("K132").Value = 10238.6001
$ws.Range("M132").Value = -7708.6001
$ws.Range("H137").Value = 12402
$ws.Range("I137").Value = 5743.6523
$ws.Range("J137").Value = 20059.1
$ws.Range("K137").Value = 17230.9569
$ws.Range("L137").Value = 60177.3
$ws.Range("M137").Value = -14680.9569
$ws.Range("N137").Value = -65277.3
$ws.Range("H138").Value = 4772.151
$ws.Range("I138").Value = 3714.2942
$ws.Range("J138").Value = 5271.6943
$ws.Range("K138").Value = 11142.8826
$ws.Range("L138").Value = 15815.0829
$ws.Range("M138").Value = -6002.882599999999
$ws.Range("N138").Value = -26095.0829

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3480.8809
$ws.Range("I32").Value = 2018.9429
$ws.Range("J32").Value = 10790.571
$ws.Range("K32").Value = 2018.9429
$ws.Range("L32").Value = 10790.571
$ws.Range("M32").Value = -1731.9429
$ws.Range("N32").Value = -11364.571
$ws.Range("H35").Value = 1557.1428
$ws.Range("I35").Value = 1557.1428
$ws.Range("K35").Value = 1557.1428
$ws.Range("M35").Value = -1151.1428
$ws.Range("H61").Value = 7433.871
$ws.Range("I61").Value = 5694.6924
$ws.Range("J61").Value = 16477.6
$ws.Range("K61").Value = 5694.6924
$ws.Range("L61").Value = 16477.6
$ws.Range("M61").Value = -5482.6924
$ws.Range("N61").Value = -16901.6
$ws.Range("H102").Value = 991
$ws.Range("I102").Value = 968.9583
$ws.Range("J102").Value = 1255.5
$ws.Range("K102").Value = 968.9583
$ws.Range("L102").Value = 1255.5
$ws.Range("M102").Value = 653.0417
$ws.Range("N102").Value = -4499.5
$ws.Range("H110").Value = 1779.2858
$ws.Range("I110").Value = 1575.1111
$ws.Range("K110").Value = 1575.1111
$ws.Range("M110").Value = 469.8888999999999
$ws.Range("H122").Value = 4529.3335
$ws.Range("I122").Value = 4529.3335
$ws.Range("K122").Value = 13588.0005
$ws.Range("M122").Value = -11138.0005
$ws.Range("H136").Value = 7433.871
$ws.Range("I136").Value = 5694.6924
$ws.Range("J136").Value = 16477.6
$ws.Range("K136").Value = 17084.0772
$ws.Range("L136").Value = 49432.8
$ws.Range("M136").Value = -14534.0772
$ws.Range("N136").Value = -54532.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").Value = ""
$ws.Range("H96").Value = 12464
$ws.Range("I96").Value = 12464
$ws.Range("K96").Value = 12464
$ws.Range("M96").Value = -9718

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 206.71428
$ws.Range("I7").Value = 203.875
$ws.Range("J7").Value = 215.8
$ws.Range("K7").Value = 203.875
$ws.Range("L7").Value = 215.8
$ws.Range("M7").Value = -90.875
$ws.Range("N7").Value = -441.8
$ws.Range("H31").Value = 6301.565
$ws.Range("I31").Value = 1492.0741
$ws.Range("J31").Value = 8540.466
$ws.Range("K31").Value = 1492.0741
$ws.Range("L31").Value = 8540.466
$ws.Range("M31").Value = -1197.0741
$ws.Range("N31").Value = -9130.466
$ws.Range("H34").Value = 6301.565
$ws.Range("I34").Value = 1492.0741
$ws.Range("J34").Value = 8540.466
$ws.Range("K34").Value = 1492.0741
$ws.Range("L34").Value = 8540.466
$ws.Range("M34").Value = -1290.0741
$ws.Range("N34").Value = -8944.466
$ws.Range("H53").Value = 74996
$ws.Range("J53").Value = 74996
$ws.Range("L53").Value = 74996
$ws.Range("N53").Value = -76210
$ws.Range("H86").Value = 5295625.5
$ws.Range("I86").Value = 7411921.5
$ws.Range("J86").Value = 4885
$ws.Range("K86").Value = 7411921.5
$ws.Range("L86").Value = 4885
$ws.Range("M86").Value = -7410798.5
$ws.Range("N86").Value = -7131
$ws.Range("H89").Value = 5295625.5
$ws.Range("I89").Value = 7411921.5
$ws.Range("J89").Value = 4885
$ws.Range("K89").Value = 37059607.5
$ws.Range("L89").Value = 24425
$ws.Range("M89").Value = -37053991.5
$ws.Range("N89").Value = -35657
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 389.72223
$ws.Range("I2").Value = 151.66667
$ws.Range("J2").Value = 627.7778
$ws.Range("K2").Value = 910.0000200000001
$ws.Range("L2").Value = 3766.6668
$ws.Range("M2").Value = -797.0000200000001
$ws.Range("N2").Value = -3992.6668
$ws.Range("H5").Value = 2521.5
$ws.Range("I5").Value = 602.4286
$ws.Range("J5").Value = 6999.3335
$ws.Range("K5").Value = 1807.2858
$ws.Range("L5").Value = 20998.0005
$ws.Range("M5").Value = -1695.2858
$ws.Range("N5").Value = -21222.0005
$ws.Range("H12").Value = 75.25
$ws.Range("I12").Value = 56.2
$ws.Range("J12").Value = 80.26316
$ws.Range("K12").Value = 168.6
$ws.Range("L12").Value = 240.78948
$ws.Range("M12").Value = 4.399999999999977
$ws.Range("N12").Value = -586.78948
$ws.Range("H38").Value = 2860.3333
$ws.Range("I38").Value = 476.1111
$ws.Range("K38").Value = 1428.3333
$ws.Range("M38").Value = -1081.3333
$ws.Range("H50").Value = 1345.711
$ws.Range("I50").Value = 399.5
$ws.Range("J50").Value = 1438.0244
$ws.Range("K50").Value = 1198.5
$ws.Range("L50").Value = 4314.0732
$ws.Range("M50").Value = -717.5
$ws.Range("N50").Value = -5276.0732
$ws.Range("H53").Value = 1345.711
$ws.Range("I53").Value = 399.5
$ws.Range("J53").Value = 1438.0244
$ws.Range("K53").Value = 1198.5
$ws.Range("L53").Value = 4314.0732
$ws.Range("M53").Value = -717.5
$ws.Range("N53").Value = -5276.0732
$ws.Range("H68").Value = 2950.261
$ws.Range("J68").Value = 3097.3572
$ws.Range("L68").Value = 9292.071599999999
$ws.Range("N68").Value = -10914.0716
$ws.Range("H71").Value = 2950.261
$ws.Range("J71").Value = 3097.3572
$ws.Range("L71").Value = 27876.2148
$ws.Range("N71").Value = -35988.2148
$ws.Range("H107").Value = 2356.25
$ws.Range("I107").Value = 2152.3635
$ws.Range("J107").Value = 4599
$ws.Range("K107").Value = 6457.0905
$ws.Range("L107").Value = 13797
$ws.Range("M107").Value = -4537.0905
$ws.Range("N107").Value = -17637
$ws.Range("H135").Value = 2521.5
$ws.Range("I135").Value = 602.4286
$ws.Range("J135").Value = 6999.3335
$ws.Range("K135").Value = 5421.8574
$ws.Range("L135").Value = 62994.0015
$ws.Range("M135").Value = -2886.8574
$ws.Range("N135").Value = -68064.0015

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 99999.5
$ws.Range("I18").Value = 99999
$ws.Range("K18").Value = 99999
$ws.Range("M18").Value = -99706
$ws.Range("H58").Value = 27377.5
$ws.Range("I58").Value = 20041
$ws.Range("J58").Value = 28844.8
$ws.Range("K58").Value = 20041
$ws.Range("L58").Value = 28844.8
$ws.Range("M58").Value = -19764
$ws.Range("N58").Value = -29398.8
$ws.Range("H122").Value = 4642.533
$ws.Range("I122").Value = 4576.5
$ws.Range("J122").Value = 4718
$ws.Range("K122").Value = 13729.5
$ws.Range("L122").Value = 14154
$ws.Range("M122").Value = -11279.5
$ws.Range("N122").Value = -19054

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1932.2307
$ws.Range("I46").Value = 1290.6666
$ws.Range("J46").Value = 2482.1428
$ws.Range("K46").Value = 1290.6666
$ws.Range("L46").Value = 2482.1428
$ws.Range("M46").Value = -1102.6666
$ws.Range("N46").Value = -2858.1428
$ws.Range("H61").Value = 1820.4706
$ws.Range("I61").Value = 1868.7273
$ws.Range("K61").Value = 1868.7273
$ws.Range("M61").Value = -1666.7273
$ws.Range("H113").Value = 1820.4706
$ws.Range("I113").Value = 1868.7273
$ws.Range("K113").Value = 1868.7273
$ws.Range("M113").Value = 301.2727
$ws.Range("H122").Value = 3469.9
$ws.Range("I122").Value = 2514.8572
$ws.Range("K122").Value = 7544.571599999999
$ws.Range("M122").Value = -5094.571599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1207.4736
$ws.Range("I107").Value = 1038.9333
$ws.Range("J107").Value = 1839.5
$ws.Range("K107").Value = 3116.7999
$ws.Range("L107").Value = 5518.5
$ws.Range("M107").Value = -1196.7999
$ws.Range("N107").Value = -9358.5
$ws.Range("H137").Value = 64500
$ws.Range("J137").Value = 69000
$ws.Range("L137").Value = 69000
$ws.Range("N137").Value = -79200
